$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9609408974647522
$ws.Range("B1").Value = 1.803401470184326
$ws.Range("C1").Value = 7.17421817779541
$ws.Range("D1").Value = 3.159863948822021
$ws.Range("E1").Value = 1.482699275016785
